$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$ws.Range("C2").Value = "z11111"
$ws.Range("C3").Value = "Z22222"
$ws.Range("C4").Value = "z33333"
$ws.Range("C5").Value = "z44444"
$ws.Range("C6").Value = "z55555"
$ws.Range("C7").Value = "z66666"
$ws.Range("C8").Value = "z77777"
$ws.Range("C9").Value = "z88888"
$ws.Range("C10").Value = "z99999"
$ws.Range("C11").Value = "z1010101010"
